# Applies the updated cryptocurrency market data (prices / 1h volume change)
# to Sheet1 of the workbook, matching the latest scrape from GitHub Actions.
#
# Columns D (Price) contain numeric-looking text (e.g. "414.59") that must stay
# plain text, exactly like the original cells. Assigning a string starting with
# a single quote (') tells Excel to store the text as-is (quote-prefixed text)
# instead of auto-converting it to a number. In a single-quoted PowerShell string,
# two single quotes ('''') produce one literal quote character, so a cell value
# is written as '''<value>' to prepend that leading quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '66.513.63'
$ws.Range('E2').Value = '  +5.90%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.594.78'
$ws.Range('E3').Value = '  +3.56%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.30%  '

# Row 5: BNB
$ws.Range('D5').Value = '''414.59'
$ws.Range('E5').Value = '  +0.07%  '

# Row 6: Solana
$ws.Range('D6').Value = '''130.46'
$ws.Range('E6').Value = '  -0.12%  '

# Row 7: XRP
$ws.Range('E7').Value = '  +3.18%  '

# Row 8: LidoStakedEther
$ws.Range('D8').Value = '3.590.31'
$ws.Range('E8').Value = '  +3.59%  '

# Row 9: USDC
$ws.Range('E9').Value = '  +0.02%  '

# Row 10: Cardano
$ws.Range('D10').Value = '''0.770'
$ws.Range('E10').Value = '  +6.01%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.173'
$ws.Range('E11').Value = '  +15.30%  '

# Row 12: ShibaInu
$ws.Range('D12').Value = '''0.0000324'
$ws.Range('E12').Value = '  +45.74%  '

# Row 13: Avalanche
$ws.Range('D13').Value = '''42.28'
$ws.Range('E13').Value = '  -0.95%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''9.84'
$ws.Range('E14').Value = '  +2.58%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '4.157.57'
$ws.Range('E15').Value = '  +3.39%  '

# Row 16: TRON
$ws.Range('E16').Value = '  -0.40%  '

# Row 17: Chainlink
$ws.Range('D17').Value = '''20.27'
$ws.Range('E17').Value = '  -1.39%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.627.99'
$ws.Range('E18').Value = '  +4.72%  '

# Row 19: Polygon
$ws.Range('D19').Value = '''1.14'
$ws.Range('E19').Value = '  +5.58%  '

# Row 20: WrappedBTC
$ws.Range('D20').Value = '66.534.45'
$ws.Range('E20').Value = '  +6.07%  '

# Row 21: Uniswap
$ws.Range('D21').Value = '''12.26'
$ws.Range('E21').Value = '  -3.20%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''445.78'
$ws.Range('E22').Value = '  -4.59%  '

# Row 23: Litecoin
$ws.Range('D23').Value = '''88.76'
$ws.Range('E23').Value = '  -2.36%  '

# Row 24: ImmutableX
$ws.Range('E24').Value = '  -3.97%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').Value = '''13.03'
$ws.Range('E25').Value = '  -1.82%  '

# Row 26: PancakeSwap
$ws.Range('E26').Value = '  +0.52%  '

# Row 27: Filecoin
$ws.Range('D27').Value = '''9.95'
$ws.Range('E27').Value = '  -5.86%  '

# Row 28: EthereumClassic
$ws.Range('D28').Value = '''35.26'
$ws.Range('E28').Value = '  +5.34%  '

# Row 29: LEO
$ws.Range('D29').Value = '''4.86'
$ws.Range('E29').Value = '  +1.09%  '

# Row 30: Toncoin
$ws.Range('E30').Value = '  +4.01%  '

# Row 31: Cosmos
$ws.Range('D31').Value = '''12.31'
$ws.Range('E31').Value = '  +2.48%  '

# Row 32: Hedera
$ws.Range('E32').Value = '  +4.10%  '

# Row 33: RenderToken
$ws.Range('D33').Value = '''7.37'
$ws.Range('E33').Value = '  -2.41%  '

# Row 34: Kaspa
$ws.Range('E34').Value = '  -4.11%  '

# Row 35: InjectiveProtocol
$ws.Range('D35').Value = '''39.64'
$ws.Range('E35').Value = '  -3.02%  '

# Row 36: Dai
$ws.Range('E36').Value = '  +0.06%  '

# Row 37: OKB
$ws.Range('D37').Value = '''56.51'
$ws.Range('E37').Value = '  -3.86%  '

# Row 38: VeChain
$ws.Range('E38').Value = '  -0.36%  '

# Row 39: PEPE
$ws.Range('D39').Value = '0.0₃0711'
$ws.Range('E39').Value = '  +28.00%  '

# Row 40: Stellar
$ws.Range('D40').Value = '''0.146'
$ws.Range('E40').Value = '  +9.28%  '

# Row 41: FirstDigitalUSD
$ws.Range('E41').Value = '  -0.02%  '

# Row 42: Stacks
$ws.Range('D42').Value = '''2.99'
$ws.Range('E42').Value = '  -2.90%  '

# Row 43: Monero
$ws.Range('D43').Value = '''148.89'
$ws.Range('E43').Value = '  +1.01%  '

# Row 44: WEMIXToken
$ws.Range('E44').Value = '  +1.02%  '

# Row 45: LidoDAOToken
$ws.Range('D45').Value = '''3.25'
$ws.Range('E45').Value = '  -2.48%  '

# Row 46: NEARProtocol
$ws.Range('E46').Value = '  -0.76%  '

# Row 47: TheGraph
$ws.Range('E47').Value = '  -3.72%  '

# Row 48: ARBITRUM
$ws.Range('E48').Value = '  -4.99%  '

# Row 49: ThetaToken
$ws.Range('D49').Value = '''2.28'
$ws.Range('E49').Value = '  -5.10%  '

# Row 50: BitcoinSV
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').Value = '''15.59'
$ws.Range('E50').Value = '  -4.89%  '

# Row 51: Celestia
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '''113.48'
$ws.Range('E51').Value = '  +3.58%  '
